# Corrigido integração com o Chrome
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Remove the two trailing rows (Renata and Alessandra) ---
$ws.Rows.Item(6).Resize(2).Delete() | Out-Null

# --- Row heights: the old rows were 75pt tall to fit the wrapped message
# text; that text is gone now, so go back to the sheet's default height. ---
$ws.Range("A1:A5").EntireRow.AutoFit()

# --- Header row ---
$ws.Range("D1").Value = "valor"

# --- Data rows: column D becomes a numeric "valor" column instead of the
# free-text "mensagem" column. ---
$ws.Range("D2").Value = 37.5
$ws.Range("D3").Value = 122.3
$ws.Range("D4").Value = 200.2

# Row 5 used to be Viviane; she's replaced by Manuela with a new phone
# number and value.
$ws.Range("B5").Value = "Manuela"
$ws.Range("C5").Value = 5532991599484
$ws.Range("D5").Value = 200

# --- Column D formatting: drop the wrap-text used for long messages, give
# it the same width as column C, and format the numbers with one decimal
# place. ---
$ws.Range("D1:D5").WrapText = $false
$ws.Range("D1:D5").HorizontalAlignment = -4108
$ws.Range("D1:D5").VerticalAlignment = -4108
$ws.Columns.Item(4).ColumnWidth = 27.85546875
$ws.Range("D1:D5").NumberFormat = "0.0"

# --- Selection / view ---
$ws.Range("G9").Select()
